$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 116.5
$ws.Range("J55").Value = 162
$ws.Range("L55").Value = 162
$ws.Range("N55").Value = -590
$ws.Range("H113").Value = 2148.2307
$ws.Range("I113").Value = 2446.3333
$ws.Range("J113").Value = 1892.7142
$ws.Range("K113").Value = 2446.3333
$ws.Range("L113").Value = 1892.7142
$ws.Range("M113").Value = 807.6667000000002
$ws.Range("N113").Value = -8400.7142
$ws.Range("H116").Value = 3018.6
$ws.Range("J116").Value = 3018.6
$ws.Range("L116").Value = 3018.6
$ws.Range("N116").Value = -9902.6
$ws.Range("H125").Value = 4152.5
$ws.Range("I125").Value = 10032
$ws.Range("J125").Value = 3312.5715
$ws.Range("K125").Value = 90288
$ws.Range("L125").Value = 29813.1435
$ws.Range("M125").Value = -87828
$ws.Range("N125").Value = -34733.1435
$ws.Range("H132").Value = 7149674.5
$ws.Range("I132").Value = 8340827
$ws.Range("K132").Value = 25022481
$ws.Range("M132").Value = -25019951
$ws.Range("H137").Value = 1403.6786
$ws.Range("J137").Value = 1122.2
$ws.Range("L137").Value = 3366.6
$ws.Range("N137").Value = -8466.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 168796
$ws.Range("I45").Value = 251449.75
$ws.Range("J45").Value = 3488.5
$ws.Range("K45").Value = 251449.75
$ws.Range("L45").Value = 3488.5
$ws.Range("M45").Value = -251072.75
$ws.Range("N45").Value = -4242.5
$ws.Range("H102").Value = 127185.125
$ws.Range("I102").Value = 252395
$ws.Range("J102").Value = 1975.25
$ws.Range("K102").Value = 252395
$ws.Range("L102").Value = 1975.25
$ws.Range("M102").Value = -250773
$ws.Range("N102").Value = -5219.25
$ws.Range("H117").Value = 39860
$ws.Range("J117").Value = 39860
$ws.Range("L117").Value = 39860
$ws.Range("N117").Value = -49038
$ws.Range("H122").Value = 2168.7334
$ws.Range("I122").Value = 2024.4546
$ws.Range("J122").Value = 2565.5
$ws.Range("K122").Value = 6073.3638
$ws.Range("L122").Value = 7696.5
$ws.Range("M122").Value = -3623.3638
$ws.Range("N122").Value = -12596.5
$ws.Range("H132").Value = 1506
$ws.Range("I132").Value = 1506.3414
$ws.Range("J132").Value = 1499
$ws.Range("K132").Value = 4519.0242
$ws.Range("L132").Value = 4497
$ws.Range("M132").Value = -1989.0242
$ws.Range("N132").Value = -9557

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H93").Value = 54866.668
$ws.Range("J93").Value = 54866.668
$ws.Range("L93").Value = 54866.668
$ws.Range("N93").Value = -58610.668
$ws.Range("H94").Value = 624
$ws.Range("I94").Value = 487.8
$ws.Range("J94").Value = 818.5714
$ws.Range("K94").Value = 487.8
$ws.Range("L94").Value = 818.5714
$ws.Range("M94").Value = -36.80000000000001
$ws.Range("N94").Value = -1720.5714
$ws.Range("H96").Value = 11400
$ws.Range("I96").Value = 6750
$ws.Range("K96").Value = 6750
$ws.Range("M96").Value = -4004
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H99").Value = 1345.2051
$ws.Range("I99").Value = 947.1539
$ws.Range("J99").Value = 1544.2307
$ws.Range("K99").Value = 947.1539
$ws.Range("L99").Value = 1544.2307
$ws.Range("M99").Value = 550.8461
$ws.Range("N99").Value = -4540.2307
$ws.Range("H101").Value = 28000
$ws.Range("J101").Value = 28000
$ws.Range("L101").Value = 28000
$ws.Range("N101").Value = -34490
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H105").Value = 145284.42
$ws.Range("I105").Value = 202576
$ws.Range("J105").Value = 2055.5
$ws.Range("K105").Value = 202576
$ws.Range("L105").Value = 2055.5
$ws.Range("M105").Value = -200829
$ws.Range("N105").Value = -5549.5
$ws.Range("H106").Value = 32000
$ws.Range("J106").Value = 32000
$ws.Range("L106").Value = 32000
$ws.Range("N106").Value = -34524
$ws.Range("H107").Value = 76959800
$ws.Range("I107").Value = 142923180
$ws.Range("J107").Value = 2525.6667
$ws.Range("K107").Value = 142923180
$ws.Range("L107").Value = 2525.6667
$ws.Range("M107").Value = -142921260
$ws.Range("N107").Value = -6365.6667
$ws.Range("H134").Value = 1870.2683
$ws.Range("I134").Value = 1612.2894
$ws.Range("J134").Value = 5138
$ws.Range("K134").Value = 4836.8682
$ws.Range("L134").Value = 15414
$ws.Range("M134").Value = -2301.8682
$ws.Range("N134").Value = -20484

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 20638.625
$ws.Range("J43").Value = 20638.625
$ws.Range("L43").Value = 20638.625
$ws.Range("N43").Value = -21006.625
$ws.Range("H51").Value = 7971.143
$ws.Range("J51").Value = 7971.143
$ws.Range("L51").Value = 7971.143
$ws.Range("N51").Value = -9443.143
$ws.Range("H61").Value = 7971.143
$ws.Range("J61").Value = 7971.143
$ws.Range("L61").Value = 7971.143
$ws.Range("N61").Value = -8667.143
$ws.Range("H68").Value = 18527.572
$ws.Range("J68").Value = 18527.572
$ws.Range("L68").Value = 18527.572
$ws.Range("N68").Value = -20025.572
$ws.Range("H71").Value = 18527.572
$ws.Range("J71").Value = 18527.572
$ws.Range("L71").Value = 55582.716
$ws.Range("N71").Value = -63070.716
$ws.Range("H74").Value = 26896.285
$ws.Range("J74").Value = 26896.285
$ws.Range("L74").Value = 26896.285
$ws.Range("N74").Value = -28644.285
$ws.Range("H77").Value = 26896.285
$ws.Range("J77").Value = 26896.285
$ws.Range("L77").Value = 80688.855
$ws.Range("N77").Value = -89424.855
$ws.Range("H101").Value = 20638.625
$ws.Range("J101").Value = 20638.625
$ws.Range("L101").Value = 20638.625
$ws.Range("N101").Value = -27128.625
$ws.Range("H107").Value = 567.80646
$ws.Range("I107").Value = 553.5909
$ws.Range("J107").Value = 602.55554
$ws.Range("K107").Value = 553.5909
$ws.Range("L107").Value = 602.55554
$ws.Range("M107").Value = 1366.4091
$ws.Range("N107").Value = -4442.55554

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 189
$ws.Range("I8").Value = 189
$ws.Range("K8").Value = 567
$ws.Range("M8").Value = -428
$ws.Range("H39").Value = 1914.091
$ws.Range("J39").Value = 2055.5
$ws.Range("L39").Value = 6166.5
$ws.Range("N39").Value = -6754.5
$ws.Range("H131").Value = 764.8
$ws.Range("I131").Value = 273.16666
$ws.Range("J131").Value = 796.18085
$ws.Range("K131").Value = 819.4999799999999
$ws.Range("L131").Value = 2388.54255
$ws.Range("M131").Value = 4220.50002
$ws.Range("N131").Value = -12468.54255
$ws.Range("H132").Value = 1667.1052
$ws.Range("J132").Value = 1898.3334
$ws.Range("L132").Value = 17085.0006
$ws.Range("N132").Value = -22145.0006
$ws.Range("H139").Value = 1967.3928
$ws.Range("I139").Value = 844.9286
$ws.Range("J139").Value = 3089.8572
$ws.Range("K139").Value = 2534.7858
$ws.Range("L139").Value = 9269.571599999999
$ws.Range("M139").Value = 2605.2142
$ws.Range("N139").Value = -19549.5716

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2611.2964
$ws.Range("I102").Value = 1881.8823
$ws.Range("K102").Value = 1881.8823
$ws.Range("M102").Value = -259.8823
$ws.Range("H107").Value = 553.1
$ws.Range("J107").Value = 456.6
$ws.Range("L107").Value = 456.6
$ws.Range("N107").Value = -4296.6
$ws.Range("H132").Value = 2374.7112
$ws.Range("I132").Value = 2273.756
$ws.Range("J132").Value = 3409.5
$ws.Range("K132").Value = 6821.268
$ws.Range("L132").Value = 10228.5
$ws.Range("M132").Value = -4291.268
$ws.Range("N132").Value = -15288.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1695
$ws.Range("I61").Value = 1781.1
$ws.Range("J61").Value = 1608.9
$ws.Range("K61").Value = 1781.1
$ws.Range("L61").Value = 1608.9
$ws.Range("M61").Value = -1579.1
$ws.Range("N61").Value = -2012.9
$ws.Range("H103").Value = 34500
$ws.Range("J103").Value = 34500
$ws.Range("L103").Value = 34500
$ws.Range("N103").Value = -36844
$ws.Range("H113").Value = 1695
$ws.Range("I113").Value = 1781.1
$ws.Range("J113").Value = 1608.9
$ws.Range("K113").Value = 1781.1
$ws.Range("L113").Value = 1608.9
$ws.Range("M113").Value = 388.9000000000001
$ws.Range("N113").Value = -5948.9

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1638.7646
$ws.Range("I136").Value = 600.56525
$ws.Range("J136").Value = 3809.5454
$ws.Range("K136").Value = 1801.69575
$ws.Range("L136").Value = 11428.6362
$ws.Range("M136").Value = 748.3042500000001
$ws.Range("N136").Value = -16528.6362
